# Fix the stray tab character in the "id" header cell (A1) and update the
# active selection to reflect the last cell touched by this edit (D1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A1 currently holds "id`t" (id followed by a tab) - replace with a clean "id"
$ws.Range("A1").Value = "id"

# Move / record the active selection on D1, matching the saved view state
$ws.Range("D1").Select()
